# Fill in missing "Source" and "Link" (JLCPCB part detail) data in the BoM sheet,
# matching the upstream KiBot-generated BoM refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BoM")

# Row 10 - D101..D113 (LED, Hubei KENTO Elec)
$ws.Range("I10").Value = "https://jlcpcb.com/partdetail/Hubei_KENTOElec-KT0603R/C2286"

# Row 12 - D201 (SK6812MINI)
$ws.Range("G12").Value = "JLCPCB"
$ws.Range("I12").Value = "https://jlcpcb.com/partdetail/OPSCOOptoelectronics-SK6812MINIHS/C2922787"
$ws.Rows(12).RowHeight = 30

# Row 15 - J201 (USB_B_Micro)
$ws.Range("G15").Value = "JLCPCB"
$ws.Range("I15").Value = "https://jlcpcb.com/partdetail/AmphenolICC-101035940001LF/C428495"
$ws.Rows(15).RowHeight = 30

# Row 16 - P8 P9 (BeagleBone_Black_Header)
$ws.Range("G16").Value = "TME"

# Row 17 - R101..R117 (UNI-ROYAL resistor)
$ws.Range("I17").Value = "https://jlcpcb.com/partdetail/12256-0402WGF1001TCE/C11702"

# Row 19 - U103 (AT24CS64-SSHM, Microchip Tech)
$ws.Range("I19").Value = "https://jlcpcb.com/partdetail/MicrochipTech-AT24C256C_SSHLT/C6482"

# Row 20 - U202 (ESP32-S3-WROOM-1)
$ws.Range("G20").Value = "JLCPCB"
$ws.Range("I20").Value = "https://jlcpcb.com/partdetail/3522416-ESP32_S3_WROOM_1UN16R8/C3013946"

# Row 21 - U101 (SP3485EN, MaxLinear)
$ws.Range("I21").Value = "https://jlcpcb.com/partdetail/MaxLinear-SP3485EN_LTR/C8963"

# Row 22 - U201 (USBLC6-2SC6)
$ws.Range("G22").Value = "JLCPCB"
$ws.Range("I22").Value = "https://jlcpcb.com/partdetail/2790619-USBLC62SC6/C2687116"
